$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '58.392.41'
$ws.Cells.Item(2, 5).Value = '  -3.35%  '
$ws.Cells.Item(3, 4).Value = '2.697.41'
$ws.Cells.Item(3, 5).Value = '  -6.82%  '
$ws.Cells.Item(4, 5).Value = '  +0.03%  '
$ws.Cells.Item(5, 4).Value = "'501.49"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -4.97%  '
$ws.Cells.Item(6, 4).Value = "'139.67"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -2.32%  '
$ws.Cells.Item(7, 4).Value = "'0.999"
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  -0.10%  '
$ws.Cells.Item(8, 4).Value = "'0.526"
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  -5.04%  '
$ws.Cells.Item(9, 4).Value = '2.705.88'
$ws.Cells.Item(9, 5).Value = '  -6.65%  '
$ws.Cells.Item(10, 2).Value = 'Toncoin'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(10, 4).Value = "'6.04"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +1.64%  '
$ws.Cells.Item(11, 2).Value = 'Dogecoin'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(11, 4).Value = "'0.104"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -3.35%  '
$ws.Cells.Item(12, 4).Value = "'0.347"
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  -3.97%  '
$ws.Cells.Item(13, 5).Value = '  +1.01%  '
$ws.Cells.Item(14, 4).Value = '3.172.72'
$ws.Cells.Item(14, 5).Value = '  -6.72%  '
$ws.Cells.Item(15, 4).Value = '58.463.21'
$ws.Cells.Item(15, 5).Value = '  -3.24%  '
$ws.Cells.Item(16, 4).Value = "'21.43"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  -5.32%  '
$ws.Cells.Item(17, 4).Value = '2.710.63'
$ws.Cells.Item(17, 5).Value = '  -6.30%  '
$ws.Cells.Item(18, 4).Value = "'0.0000134"
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  -5.63%  '
$ws.Cells.Item(19, 4).Value = "'4.72"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -6.38%  '
$ws.Cells.Item(20, 4).Value = "'10.85"
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -7.09%  '
$ws.Cells.Item(21, 4).Value = "'333.56"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -8.37%  '
$ws.Cells.Item(22, 4).Value = "'6.20"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -5.96%  '
$ws.Cells.Item(23, 4).Value = "'0.996"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -0.34%  '
$ws.Cells.Item(24, 4).Value = "'5.63"
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -0.14%  '
$ws.Cells.Item(25, 4).Value = "'62.94"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -2.03%  '
$ws.Cells.Item(26, 4).Value = "'0.173"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -2.93%  '
$ws.Cells.Item(27, 4).Value = "'0.424"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -6.20%  '
$ws.Cells.Item(28, 4).Value = "'0.996"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -0.23%  '
$ws.Cells.Item(29, 4).Value = "'7.40"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -5.25%  '
$ws.Cells.Item(30, 4).Value = '0.0₃0817'
$ws.Cells.Item(30, 5).Value = '  -5.14%  '
$ws.Cells.Item(31, 5).Value = '  -0.10%  '
$ws.Cells.Item(32, 2).Value = 'PancakeSwap'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(32, 4).Value = "'1.59"
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -4.69%  '
$ws.Cells.Item(33, 2).Value = 'EthereumClassic'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(33, 4).Value = "'19.06"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -2.90%  '
$ws.Cells.Item(34, 4).Value = "'151.23"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +6.02%  '
$ws.Cells.Item(35, 4).Value = "'5.37"
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  -3.60%  '
$ws.Cells.Item(36, 4).Value = "'4.15"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -4.64%  '
$ws.Cells.Item(37, 4).Value = "'0.928"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -7.05%  '
$ws.Cells.Item(38, 4).Value = "'1.12"
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -7.04%  '
$ws.Cells.Item(39, 4).Value = "'35.24"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -6.67%  '
$ws.Cells.Item(40, 4).Value = "'3.55"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -3.24%  '
$ws.Cells.Item(41, 4).Value = "'1.37"
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -8.01%  '
$ws.Cells.Item(42, 4).Value = '2.177.69'
$ws.Cells.Item(42, 5).Value = '  -6.29%  '
$ws.Cells.Item(43, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(43, 4).Value = "'0.997"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -0.08%  '
$ws.Cells.Item(44, 2).Value = 'Hedera'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(44, 4).Value = "'0.0556"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -3.58%  '
$ws.Cells.Item(45, 4).Value = "'0.593"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -7.83%  '
$ws.Cells.Item(46, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Cells.Item(46, 4).Value = "'10.36"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +0.38%  '
$ws.Cells.Item(47, 2).Value = 'EnergySwap'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(47, 4).Value = "'18.75"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -9.08%  '
$ws.Cells.Item(48, 4).Value = "'4.64"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -5.72%  '
$ws.Cells.Item(49, 4).Value = "'0.0225"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -3.54%  '
$ws.Cells.Item(50, 4).Value = "'0.0885"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -5.23%  '
$ws.Cells.Item(51, 4).Value = "'17.84"
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -3.16%  '
